# The commit removes the "4.2" backlog item (public course selection story)
# from the Sheet1 backlog table. In the original sheet this lived on row 17
# (Story ID "4.2" / "Public course selection"), with rows 18-20 ("4.3",
# "4.4", "4.5") following it. Deleting the entire row shifts the remaining
# rows up by one, matching the data seen in the target workbook (dimension
# shrinks from A1:I37 to A1:I36, and the corresponding shared strings for
# the removed story are no longer referenced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Delete()

# Reflect the author's final cursor position/selection in the saved view.
$ws.Range("B23").Select() | Out-Null
